# edit.ps1
# Applies the "Mejora del informe ejecutivo" change:
# Inserts a new "Entendimiento de los datos:" section right after the
# "Informe Ejecutivo" title, before the existing "Hallazgos del Modelo de
# Regresion:" section.

$d = $word.ActiveDocument

$textHeader = 'Entendimiento de los datos:'
$textD = 'Al explorar y analizar los datos, se incorporaron dos variables adicionales significativas: "Región" (region) y "Nivel de Ingreso" (income). Estas variables resultaron ser fundamentales en el modelo de regresión, desempeñando un papel crucial en la determinación del Producto Interno Bruto (PIB) de los países. La variable "Región" reveló las disparidades económicas y de desarrollo entre diferentes partes del mundo, permitiendo una comprensión más profunda de los patrones geográficos en el crecimiento económico. Por otro lado, la categorización del ingreso proporcionó una visión clara de la capacidad económica de cada país, distinguiendo entre economías de bajos, medianos y altos ingresos.'
$f1 = 'La inclusión y validación de estas variables en el modelo resaltan la importancia de considerar no solo los indicadores de salud y desarrollo interno, sino también los factores externos y estructurales que influyen en la economía de un país. La '
$f2 = 'r'
$f3 = 'egión ofrece información sobre contextos históricos y culturales únicos que impactan las economías locales, mientras que el '
$f4 = 'n'
$f5 = 'ivel de '
$f6 = 'i'
$f7 = 'ngreso proporciona una instantánea de la capacidad financiera de un país en el escenario mundial.  Este análisis reafirma la interconexión compleja entre factores geográficos, demográficos y económicos en la determinación del desarrollo económico de un país. Al entender y reconocer la influencia de estas variables adicionales, las políticas públicas pueden ser diseñadas de manera más precisa y eficaz, adaptándose a las necesidades específicas de cada región y nivel de ingreso para fomentar un crecimiento económico sostenible y equitativo.'

# Combine the text for paragraph F (runs r2..r8 in the diff all share the
# same run formatting, so they collapse naturally into one logical run of
# text when typed/inserted together).
$textF = $f1 + $f2 + $f3 + $f4 + $f5 + $f6 + $f7

# The anchor is the (currently) empty paragraph right after the
# "Informe Ejecutivo" heading paragraph. New paragraphs are inserted
# immediately before it, pushing it further down -- it ends up directly
# above "Hallazgos del Modelo de Regresion:" again, unchanged.
$anchorPara = $d.Paragraphs.Item(2)
$anchorRange = $anchorPara.Range

$combined = "`r" + $textHeader + "`r`r" + $textD + "`r`r" + $textF + "`r"
$anchorRange.InsertBefore($combined)

# Paragraphs 2-7 are now the new ones:
#  2: empty (bold heading style)
#  3: "Entendimiento de los datos:" (bold heading style)
#  4: empty (bold heading style)
#  5: paragraph D (body text, not bold, Spanish (Mexico))
#  6: empty (body text style)
#  7: paragraph F (body text, not bold, Spanish (Mexico))
for ($i = 5; $i -le 7; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pr = $p.Range
    $pr.Font.Bold = 0
    $pr.Font.LanguageID = "es-MX"
}

Write-Host "Inserted 'Entendimiento de los datos' section."
